# Fruta / hortaliza, semanal
# Insert two new weekly data rows (before the current row 462) for
# "Poroto verde" at Mercado Mayorista Lo Valledor de Santiago, origin Peru,
# date 2021-09-09 (serial 44448), and push the existing rows 462:471 down
# to 464:473.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 462, shifting rows 462:471 down to 464:473.
$ws.Range("A462:A463").EntireRow.Insert()

# New row 462: Magnum / Primera, origin Peru
$ws.Cells.Item(462, 1).Value = 6
$ws.Cells.Item(462, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(462, 3).Value = "Metropolitana"
$ws.Cells.Item(462, 4).Value = 44448
$ws.Cells.Item(462, 4).NumberFormat = $ws.Cells.Item(466, 4).NumberFormat
$ws.Cells.Item(462, 5).Value = 13
$ws.Cells.Item(462, 6).Value = 100112031
$ws.Cells.Item(462, 7).Value = "Poroto verde"
$ws.Cells.Item(462, 8).Value = "Magnum"
$ws.Cells.Item(462, 9).Value = "Primera"
$ws.Cells.Item(462, 10).Value = 400
$ws.Cells.Item(462, 11).Value = 37000
$ws.Cells.Item(462, 12).Value = 38000
$ws.Cells.Item(462, 13).Value = 37575
$ws.Cells.Item(462, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(462, 15).Value = "Perú"
$ws.Cells.Item(462, 16).Value = 1503
$ws.Cells.Item(462, 17).Value = 25
$ws.Cells.Item(462, 18).Value = "Hortaliza"

# New row 463: Sin especificar / Primera, origin Peru
$ws.Cells.Item(463, 1).Value = 6
$ws.Cells.Item(463, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(463, 3).Value = "Metropolitana"
$ws.Cells.Item(463, 4).Value = 44448
$ws.Cells.Item(463, 4).NumberFormat = $ws.Cells.Item(466, 4).NumberFormat
$ws.Cells.Item(463, 5).Value = 13
$ws.Cells.Item(463, 6).Value = 100112031
$ws.Cells.Item(463, 7).Value = "Poroto verde"
$ws.Cells.Item(463, 8).Value = "Sin especificar"
$ws.Cells.Item(463, 9).Value = "Primera"
$ws.Cells.Item(463, 10).Value = 400
$ws.Cells.Item(463, 11).Value = 35000
$ws.Cells.Item(463, 12).Value = 36000
$ws.Cells.Item(463, 13).Value = 35575
$ws.Cells.Item(463, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(463, 15).Value = "Perú"
$ws.Cells.Item(463, 16).Value = 1423
$ws.Cells.Item(463, 17).Value = 25
$ws.Cells.Item(463, 18).Value = "Hortaliza"
